$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.227.31'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '1.583.19'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '1.805.73'
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").Value = '1.580.68'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '26.246.71'
$ws.Range("E17").Value = '  -1.64%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '206.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("E23").Value = '  -3.49%  '
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '1.286.68'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.52%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.768'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").Value = '1.718.63'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0509'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  +0.08%  '
